$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new blank column before the old "E" column (Type 2 block).
#    This pushes E..I -> F..J and automatically creates two new blank
#    columns: the new E (between the 'Type 1' and 'Type 2' blocks) and
#    the new I (between the 'Type 2' block and the 'Avg' column),
#    since the sheet had no data in the old H column.
# ---------------------------------------------------------------------
$ws.Range("E1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2. Rename the header labels.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Type 1 ('Cells')"
$ws.Range("F1").Value = "Type 2 ('Debris')"

# ---------------------------------------------------------------------
# 2b. Re-enter the "D" column formula cell-by-cell (this breaks the
#     shared formula group that D inherited from the original file).
# ---------------------------------------------------------------------
$ws.Range("D2").Formula = "=(C2+B2)/2*10^4"
$ws.Range("D3").Formula = "=(C3+B3)/2*10^4"
$ws.Range("D4").Formula = "=(C4+B4)/2*10^4"
$ws.Range("D5").Formula = "=(C5+B5)/2*10^4"
$ws.Range("D6").Formula = "=(C6+B6)/2*10^4"
$ws.Range("D7").Formula = "=(C7+B7)/2*10^4"
$ws.Range("D8").Formula = "=(C8+B8)/2*10^4"
$ws.Range("D9").Formula = "=(C9+B9)/2*10^4"

# ---------------------------------------------------------------------
# 3. Fill the new "E" column (half of D) for rows 2-9.
# ---------------------------------------------------------------------
$ws.Range("E2").Formula = "=D2/2"
$ws.Range("E3:E9").Formula = "=D3/2"

# ---------------------------------------------------------------------
# 4. Fill the new "I" column (half of H) for rows 2-9.
# ---------------------------------------------------------------------
$ws.Range("I2").Formula = "=H2/2"
$ws.Range("I3:I9").Formula = "=H3/2"
$ws.Range("I2:I9").Font.Color = 0

# ---------------------------------------------------------------------
# 5. Add the halved B/C and F/G value rows (11-18), plus blank styled
#    D/E cells through row 19.
# ---------------------------------------------------------------------
$ws.Range("B11").Formula = "=B2/2"
$ws.Range("C11").Formula = "=C2/2"
$ws.Range("F11").Formula = "=F2/2"
$ws.Range("G11").Formula = "=G2/2"

$ws.Range("B12:C19").Formula = "=B3/2"
$ws.Range("B19:C19").ClearContents()

$ws.Range("F12:G12").Formula = "=F3/2"
$ws.Range("F13:G13").Formula = "=F4/2"
$ws.Range("F14:G14").Formula = "=F5/2"
$ws.Range("F15:G15").Formula = "=F6/2"
$ws.Range("F16:G16").Formula = "=F7/2"
$ws.Range("F17:G17").Formula = "=F8/2"
$ws.Range("F18:G18").Formula = "=F9/2"

# ---------------------------------------------------------------------
# 6. D/E carry the "black font" style all the way to row 19 (style was
#    inherited down the column from the original D column formatting).
# ---------------------------------------------------------------------
$ws.Range("D11:E19").Font.Color = 0

[void]$ws.Range("E22").Select()
$excel.ActiveWindow.Zoom = 98

Write-Host "done"
